$wb = $excel.ActiveWorkbook

# ---- Blackbox sheet (test case names for BB_sortByLimitingFactor) ----
$wsB = $wb.Worksheets.Item("Blackbox")
$wsB.Range("A3").Value = "BB_sortByLimitingFactor1"
$wsB.Range("A4").Value = "BB_sortByLimitingFactor2"
$wsB.Range("A5").Value = "BB_sortByLimitingFactor3"
$wsB.Range("A6").Value = "BB_sortByLimitingFactor4"
$wsB.Range("A7").Value = "BB_sortByLimitingFactor5"

# widen column A so the longer test case names are readable
$wsB.Columns.Item(1).ColumnWidth = 23.16

# ---- Whitebox sheet (test case names for WB_SortByLimitingFactor) ----
$wsW = $wb.Worksheets.Item("Whitebox")
$wsW.Range("A3").Value = "WB_SortByLimitingFactor1"
$wsW.Range("A4").Value = "WB_SortByLimitingFactor2"
$wsW.Range("A5").Value = "WB_SortByLimitingFactor3"
$wsW.Range("A6").Value = "WB_SortByLimitingFactor4"

# widen column A so the longer test case names are readable
$wsW.Columns.Item(1).ColumnWidth = 23.66

# ---- update selections / view state (Blackbox first, Whitebox last so the
# Whitebox tab remains the active tab, matching the original workbook) ----
$wsB.Range("C9").Select()
$wsW.Range("A6").Select()
